# Update "想去人数" (F column) values across the workbook's sheets.
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 368
$ws1.Range("F4").Value  = 143
$ws1.Range("F5").Value  = 1270
$ws1.Range("F7").Value  = 2438
$ws1.Range("F8").Value  = 864
$ws1.Range("F9").Value  = 18426
$ws1.Range("F10").Value = 48
$ws1.Range("F11").Value = 1851
$ws1.Range("F12").Value = 650
$ws1.Range("F13").Value = 592
$ws1.Range("F14").Value = 314
$ws1.Range("F15").Value = 588
$ws1.Range("F16").Value = 189
$ws1.Range("F17").Value = 189
$ws1.Range("F18").Value = 64
$ws1.Range("F19").Value = 313
$ws1.Range("F20").Value = 160
$ws1.Range("F21").Value = 90
$ws1.Range("F22").Value = 13

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value  = 110
$ws2.Range("F14").Value = 65

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5851
$ws3.Range("F3").Value = 546
$ws3.Range("F4").Value = 550

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 5851
$ws4.Range("F4").Value  = 546
$ws4.Range("F5").Value  = 550
$ws4.Range("F6").Value  = 368
$ws4.Range("F8").Value  = 143
$ws4.Range("F10").Value = 1270
$ws4.Range("F15").Value = 2438
$ws4.Range("F16").Value = 864
$ws4.Range("F17").Value = 18426
$ws4.Range("F18").Value = 48
$ws4.Range("F20").Value = 110
$ws4.Range("F21").Value = 111
$ws4.Range("F22").Value = 1851
$ws4.Range("F23").Value = 650
$ws4.Range("F25").Value = 592
$ws4.Range("F26").Value = 314
$ws4.Range("F27").Value = 588
$ws4.Range("F28").Value = 189
$ws4.Range("F29").Value = 189
$ws4.Range("F31").Value = 64
$ws4.Range("F34").Value = 313
$ws4.Range("F36").Value = 65
$ws4.Range("F37").Value = 160
$ws4.Range("F39").Value = 90
$ws4.Range("F43").Value = 13
